$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet="ALC"; Row=4; Cols=@{H=196.25; I=81.42857; K=81.42857; M=32.57143000000001}},
    @{Sheet="ALC"; Row=41; Cols=@{H=1028; I=1512.5; J=335.85715; K=1512.5; L=335.85715; M=-1072.5; N=-1215.85715}},
    @{Sheet="ALC"; Row=62; Cols=@{H=5094.75; I=4746.8; J=5343.2856; K=4746.8; L=5343.2856; M=-4122.8; N=-6591.2856}},
    @{Sheet="ALC"; Row=64; Cols=@{H=4260.7144; I=3700; K=3700; M=-3452}},
    @{Sheet="ALC"; Row=65; Cols=@{H=5094.75; I=4746.8; J=5343.2856; K=23734; L=26716.428; M=-20614; N=-32956.428}},
    @{Sheet="ALC"; Row=67; Cols=@{H=4260.7144; I=3700; K=3700; M=-2842}},
    @{Sheet="ALC"; Row=98; Cols=@{H=3027.9; I=2220; K=2220; M=-722}},
    @{Sheet="ALC"; Row=122; Cols=@{H=3027.9; I=2220; K=6660; M=-4210}},
    @{Sheet="ALC"; Row=132; Cols=@{H=2585.2593; I=2075.2083; J=6665.6665; K=6225.624899999999; L=19996.9995; M=-3695.624899999999; N=-25056.9995}},
    @{Sheet="ALC"; Row=138; Cols=@{H=2165.614; J=2399.7646; L=7199.293799999999; N=-17479.2938}},
    @{Sheet="ARM"; Row=32; Cols=@{H=22897.158; I=5366.1816; K=5366.1816; M=-5079.1816}},
    @{Sheet="ARM"; Row=45; Cols=@{H=5364.5884; I=6507.476; K=6507.476; M=-6130.476}},
    @{Sheet="ARM"; Row=61; Cols=@{H=2462; I=2482.25; J=2454.6365; K=2482.25; L=2454.6365; M=-2270.25; N=-2878.6365}},
    @{Sheet="ARM"; Row=74; Cols=@{H=3660.3; I=3654.3333; K=3654.3333; M=-2780.3333}},
    @{Sheet="ARM"; Row=77; Cols=@{H=3660.3; I=3654.3333; K=18271.6665; M=-13903.6665}},
    @{Sheet="ARM"; Row=118; Cols=@{H=60000; J=60000; L=60000; N=-63314}},
    @{Sheet="ARM"; Row=122; Cols=@{H=2258.5; I=2288.92; K=6866.76; M=-4416.76}},
    @{Sheet="ARM"; Row=132; Cols=@{H=2449.1052; I=2473.2778; K=7419.8334; M=-4889.8334}},
    @{Sheet="ARM"; Row=133; Cols=@{H=104629.5; J=104629.5; L=104629.5; N=-109689.5}},
    @{Sheet="ARM"; Row=136; Cols=@{H=2462; I=2482.25; J=2454.6365; K=7446.75; L=7363.9095; M=-4896.75; N=-12463.9095}},
    @{Sheet="BSM"; Row=86; Cols=@{H=1886.4762; I=1898.2667; J=1857; K=1898.2667; L=1857; M=-775.2666999999999; N=-4103}},
    @{Sheet="BSM"; Row=89; Cols=@{H=1886.4762; I=1898.2667; J=1857; K=9491.333499999999; L=9285; M=-3875.333499999999; N=-20517}},
    @{Sheet="BSM"; Row=94; Cols=@{H=1618.4807; I=1329.875; J=2580.5; K=1329.875; L=2580.5; M=-878.875; N=-3482.5}},
    @{Sheet="BSM"; Row=99; Cols=@{H=5858.7856; J=3750; L=3750; N=-6746}},
    @{Sheet="BSM"; Row=107; Cols=@{H=25142.883; I=34873.832; J=2686.8462; K=34873.832; L=2686.8462; M=-32953.832; N=-6526.8462}},
    @{Sheet="CRP"; Row=70; Cols=@{H=23090; J=23090; L=23090; N=-23720}},
    @{Sheet="CRP"; Row=73; Cols=@{H=23090; J=23090; L=23090; N=-25274}},
    @{Sheet="CRP"; Row=76; Cols=@{H=4906.75; I=4906.75; K=4906.75; M=-4591.75}},
    @{Sheet="CRP"; Row=79; Cols=@{H=4906.75; I=4906.75; K=4906.75; M=-3814.75}},
    @{Sheet="CRP"; Row=107; Cols=@{H=1731.9286; J=1842.2307; L=1842.2307; N=-5682.2307}},
    @{Sheet="CRP"; Row=134; Cols=@{H=4407.2666; I=4365; K=13095; M=-10560}},
    @{Sheet="CUL"; Row=4; Cols=@{H=13975716; I=3444034.8; K=10332104.4; M=-10331992.4}},
    @{Sheet="CUL"; Row=113; Cols=@{H=853.2069; J=875.65; L=2626.95; N=-6966.95}},
    @{Sheet="GSM"; Row=2; Cols=@{H=303.3; I=312; J=285.9; K=312; L=285.9; M=-199; N=-511.9}},
    @{Sheet="GSM"; Row=80; Cols=@{H=4546.25; J=5860.1113; L=5860.1113; N=-7856.1113}},
    @{Sheet="GSM"; Row=83; Cols=@{H=4546.25; J=5860.1113; L=29300.5565; N=-39284.5565}},
    @{Sheet="GSM"; Row=97; Cols=@{H=14989.489; J=1312.25; L=1312.25; N=-2304.25}},
    @{Sheet="GSM"; Row=102; Cols=@{H=2209.8948; I=2058.1765; K=2058.1765; M=-436.1765}},
    @{Sheet="GSM"; Row=122; Cols=@{H=1780.3334; I=1730.7142; K=5192.142599999999; M=-2742.142599999999}},
    @{Sheet="GSM"; Row=132; Cols=@{H=4252.212; I=4276.2; J=4215.3076; K=12828.6; L=12645.9228; M=-10298.6; N=-17705.9228}},
    @{Sheet="LTW"; Row=40; Cols=@{H=5156.174; I=2300.6667; K=2300.6667; M=-2164.6667}},
    @{Sheet="LTW"; Row=46; Cols=@{H=74220.5; I=423123; K=423123; M=-422935}},
    @{Sheet="LTW"; Row=82; Cols=@{H=2557.6; J=2459.8; L=2459.8; N=-3181.8}},
    @{Sheet="LTW"; Row=85; Cols=@{H=2557.6; J=2459.8; L=2459.8; N=-4955.8}},
    @{Sheet="LTW"; Row=107; Cols=@{H=3004.6667; I=3004.6667; K=3004.6667; M=-1084.6667}},
    @{Sheet="LTW"; Row=139; Cols=@{H=99888.78; J=108624.875; L=108624.875; N=-118904.875}},
    @{Sheet="WVR"; Row=81; Cols=@{H=965.625; I=965.625; J=0; K=1931.25; L=0; M=-870.25; N="__REMOVE__"}},
    @{Sheet="WVR"; Row=84; Cols=@{H=965.625; I=965.625; J=0; K=9656.25; L=0; M=-4352.25; N="__BLANK__"}},
    @{Sheet="WVR"; Row=113; Cols=@{H=545.8; I=469.26315; J=2000; K=1407.78945; L=6000; M=762.21055; N=-10340}},
    @{Sheet="WVR"; Row=132; Cols=@{H=862480.8; I=1101380.5; J=2441.9; K=3304141.5; L=7325.700000000001; M=-3301611.5; N=-12385.7}},
)

# Apply each cell update. A couple of rows (WVR 81/84) had their LeveProfitHQ (N)
# value cleared out entirely as part of the source edit, so those are special-cased
# with ClearContents() instead of a numeric assignment.
foreach ($item in $changes) {
    $ws = $wb.Worksheets.Item($item.Sheet)
    foreach ($col in $item.Cols.Keys) {
        $addr = "$col$($item.Row)"
        $val = $item.Cols[$col]
        if ($val -eq "__REMOVE__" -or $val -eq "__BLANK__") {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

Write-Host "Applied $($changes.Count) row updates"
